$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 13 and row 14 for columns A, B, D, E, F, G, H, Q, R
$row13_A = $ws.Range("A13").Value2
$row13_B = $ws.Range("B13").Value2
$row13_D = $ws.Range("D13").Value2
$row13_E = $ws.Range("E13").Value2
$row13_F = $ws.Range("F13").Value2
$row13_G = $ws.Range("G13").Value2
$row13_H = $ws.Range("H13").Value2
$row13_Q = $ws.Range("Q13").Value2
$row13_R = $ws.Range("R13").Value2

$row14_A = $ws.Range("A14").Value2
$row14_B = $ws.Range("B14").Value2
$row14_D = $ws.Range("D14").Value2
$row14_E = $ws.Range("E14").Value2
$row14_F = $ws.Range("F14").Value2
$row14_G = $ws.Range("G14").Value2
$row14_H = $ws.Range("H14").Value2
$row14_Q = $ws.Range("Q14").Value2
$row14_R = $ws.Range("R14").Value2

$ws.Range("A13").Value = $row14_A
$ws.Range("B13").Value = $row14_B
$ws.Range("D13").Value = $row14_D
$ws.Range("E13").Value = $row14_E
$ws.Range("F13").Value = $row14_F
$ws.Range("G13").Value = $row14_G
$ws.Range("H13").Value = $row14_H
$ws.Range("Q13").Value = $row14_Q
$ws.Range("R13").Value = $row14_R

$ws.Range("A14").Value = $row13_A
$ws.Range("B14").Value = $row13_B
$ws.Range("D14").Value = $row13_D
$ws.Range("E14").Value = $row13_E
$ws.Range("F14").Value = $row13_F
$ws.Range("G14").Value = $row13_G
$ws.Range("H14").Value = $row13_H
$ws.Range("Q14").Value = $row13_Q
$ws.Range("R14").Value = $row13_R
